$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.412.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.797.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.60%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3814'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3466'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.88%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.204'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07522'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.497'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.13%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.795.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.068'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001099'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06649'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '85.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.524'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.390.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.52%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.54'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.431'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.578'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.501'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.85'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.999.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.75%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.63%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.064'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.152'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08723'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.46%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.29'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.700'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.457'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6906'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.32%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.970'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.34%  '

$ws.Range("E40").Value = '  +1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06377'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02341'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.276'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.29%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6468'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.876'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.127'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.13'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07196'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.49%  '
